# To-do_list.xlsx update:
#  - Fix typo "popuplação" -> "população"
#  - Reorder a handful of to-do items (rotate B5:B9 so "Gerar o MER..." moves
#    to the bottom of that block and the rest shift up one row)
#  - Add a "Status" column (C) with two annotations
#  - Add a trailing blank styled row at the bottom of the merged priority cell
#  - Misc cosmetic tweaks (row height, selection, column width)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteValues  = -4163
$xlPasteFormats = -4122
$xlRight        = -4152
$xlCenter       = -4108

# ---------------------------------------------------------------------------
# 1. Fix the typo in the shared string used by B4 ("popuplação" -> "população")
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Uso de população passada na geração de nova chave"

# ---------------------------------------------------------------------------
# 2. Rotate B5:B9 so that "Gerar o MER do banco de dados" (currently B5) ends
#    up at B9, and the 4 items below it (B6..B9) each shift up one row.
#    Rich-text runs (B6 "Escrever o About e o Help", B9 "Adicionar ferramenta
#    de log") are moved with Copy/PasteSpecial(values) so their formatted
#    runs survive instead of being flattened to plain text.
# ---------------------------------------------------------------------------
$stagedB5 = $ws.Range("B5").Text

$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial($xlPasteValues)

$ws.Range("B7").Copy()
$ws.Range("B6").PasteSpecial($xlPasteValues)

$ws.Range("B8").Copy()
$ws.Range("B7").PasteSpecial($xlPasteValues)

$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial($xlPasteValues)

$ws.Range("B9").Value = $stagedB5
$ws.Application.CutCopyMode = $false

# B7 ("Revisar monografia" after the rotation) needs to pick up the
# orange/medium-priority formatting class (matching B8/B9) instead of the
# red/high-priority one it inherited from its old "Revisar" position.
$ws.Range("B8").Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. New "Status" column (C)
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false
$ws.Range("C1").Value = "Status"

$ws.Range("C5").Value = "TERMINADO"
$ws.Range("C4").Value = "QUASE TERMINADO"
$ws.Range("C4:C5").HorizontalAlignment = $xlRight

$ws.Columns.Item(3).ColumnWidth = 37.14

# ---------------------------------------------------------------------------
# 4. New trailing blank row (19) styled like the rest of the merged column A
#    (centered, wrapped) cell below the list
# ---------------------------------------------------------------------------
$ws.Range("A19").WrapText = $true
$ws.Range("A19").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# 5. Cosmetic tweaks
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 15
$ws.Range("D7").Select()
